$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header layout (row 1): A=ID, B=ServerID, C=Name, D=MaxOnline, E=CpuCount, F=IP, G=Port
# Fill in row 2 with the new LoginServer entry. Value-assignment order matters
# because it drives the order new entries are appended to sharedStrings.xml
# (IP first, then the repeated "LoginServer_1" name/id, then the ServerID code).
$ws.Range("F2").Value = "127.0.0.1"
$ws.Range("F2").NumberFormat = "@"

$ws.Range("A2").Value = "LoginServer_1"
$ws.Range("B2").Value = "000106001"
$ws.Range("C2").Value = "LoginServer_1"
$ws.Range("C2").NumberFormat = "@"

$ws.Range("D2").Value = 5000
$ws.Range("E2").Value = 1
$ws.Range("G2").Value = 6001

# The F column's TRUE/FALSE list validation now starts below the populated
# row 2 (it used to cover F2:F1048576, now it is F3:F1048576).
$ws.Range("F2:F1048576").Validation.Delete()
$ws.Range("F3:F1048576").Validation.Add(3, 1, 1, '"TRUE,FALSE"')

# Move the active selection to G5 (single cell), matching the saved view state.
$ws.Range("G5").Select()
